# EPBDS-9427 Fix NPE in RulesFrontend proxy
# Add a new "worldHello(Integer i, String s)" method block to the sheet,
# mirroring the existing B3:C4 "worldHello(int hour)" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20: method signature block (mirrors B3:C3 formatting - boxed
# border, centered).
$ws.Range("B20:C20").Merge()
$ws.Range("B20:C20").Borders.LineStyle = 1
$ws.Range("B20:C20").HorizontalAlignment = -4108

# New row 21: method body block (mirrors B4:C4 formatting - boxed border).
$ws.Range("B21:C21").Merge()
$ws.Range("B21:C21").Borders.LineStyle = 1

# Shared string for the "return ..." text is registered before the
# "Method ..." signature text, matching the order they were authored in.
$ws.Range("B21").Value = "return ""i: "" + I +"" s: "" + s;"
$ws.Range("B20").Value = "Method String worldHello(Integer i, String s)"

# Update selection to match target state
$ws.Range("B3:C3").Select()
